# Rename the 10 "Residential Aged Care" metric names in column C (rows 72-81)
# of the "Metrics" table, moving the "Residential Aged Care" qualifier from the
# front of the text to the end, and reshuffling a couple of word orders.
#
# Each worksheet cell keeps referring to the same underlying metric (e.g. the
# cell that used to read "# Residential Aged Care Resident Cases (Daily)"
# still refers to that same concept), only the label text itself changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newLabels = @{
    72 = "# Resident Cases (Daily) - Residential Aged Care"
    73 = "# Resident Cases (Weekly) - Residential Aged Care"
    74 = "# Staff Cases (Weekly) - Residential Aged Care"
    75 = "# Staff Cases (Weekly) per 1M -  Residential Aged Care"
    76 = "# Molnupiravir Prescriptions - Residential Aged Care"
    77 = "# Molnupiravir Prescriptions (Daily) - Residential Aged Care"
    78 = "% Molnupiravir Prescriptions per Case - Residential Aged Care"
    79 = "# Paxlovid Prescriptions - Residential Aged Care"
    80 = "# Paxlovid Prescriptions (Daily) -  Residential Aged Care"
    81 = "% Paxlovid Prescriptions per Case - Residential Aged Care"
}

foreach ($row in $newLabels.Keys) {
    $ws.Range("C$row").Value = $newLabels[$row]
}

# Update the selected cell as recorded in the saved workbook view.
$ws.Range("C82").Select()
